$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'321.84"
$ws.Range("E2").Value = "'-3.07%"
$ws.Range("E3").Value = "'-6.02%"
$ws.Range("D4").Value = "'5.211"
$ws.Range("E4").Value = "'-7.36%"
$ws.Range("D5").Value = "'0.08191"
$ws.Range("E5").Value = "'-2.03%"
$ws.Range("D6").Value = "'4.323"
$ws.Range("E6").Value = "'-2.87%"
$ws.Range("D7").Value = "'1.836"
$ws.Range("E7").Value = "'-10.07%"
$ws.Range("D8").Value = "'0.9350"
$ws.Range("E8").Value = "'-4.10%"
$ws.Range("E9").Value = "'-3.88%"
$ws.Range("E10").Value = "'-2.71%"
$ws.Range("D11").Value = "'0.09417"
$ws.Range("E11").Value = "'-5.64%"
$ws.Range("D12").Value = "'0.04621"
$ws.Range("E12").Value = "'-0.90%"
$ws.Range("D13").Value = "'7.388"
$ws.Range("E13").Value = "'-28.73%"
$ws.Range("D14").Value = "'0.1058"
$ws.Range("E14").Value = "'-0.14%"
$ws.Range("D15").Value = "'0.001299"
$ws.Range("E15").Value = "'0.43%"
$ws.Range("D16").Value = "'0.005765"
$ws.Range("E16").Value = "'-4.85%"
$ws.Range("D17").Value = "'3.357"
$ws.Range("D18").Value = "'2.511"
$ws.Range("E18").Value = "'-1.91%"
$ws.Range("D19").Value = "'0.3339"
$ws.Range("E19").Value = "'-0.75%"
$ws.Range("D20").Value = "'0.1387"
$ws.Range("E20").Value = "'-0.35%"
$ws.Range("D22").Value = "'0.04151"
$ws.Range("E22").Value = "'-1.28%"
$ws.Range("D23").Value = "'0.001250"
$ws.Range("E23").Value = "'-4.77%"
$ws.Range("D24").Value = "'0.004323"
$ws.Range("E24").Value = "'-5.98%"
$ws.Range("D25").Value = "'0.0001100"
$ws.Range("E25").Value = "'-15.50%"
$ws.Range("D26").Value = "'0.0002979"
$ws.Range("D38").Value = "'0.02723"
$ws.Range("E38").Value = "'-1.57%"
$ws.Range("D39").Value = "'0.05549"
$ws.Range("E39").Value = "'-4.44%"
$ws.Range("D40").Value = "'0.007963"
$ws.Range("E40").Value = "'2.75%"
$ws.Range("D41").Value = "'0.1396"
$ws.Range("E41").Value = "'-2.90%"
$ws.Range("D42").Value = "'0.006537"
$ws.Range("E42").Value = "'-10.06%"
$ws.Range("D43").Value = "'0.002092"
$ws.Range("E43").Value = "'-1.18%"
$ws.Range("D44").Value = "'0.007478"
$ws.Range("E44").Value = "'-7.48%"
$ws.Range("D45").Value = "'0.3204"
$ws.Range("E45").Value = "'-5.89%"
$ws.Range("D46").Value = "'0.00006966"
$ws.Range("E46").Value = "'-4.60%"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'-0.16%"
$ws.Range("D48").Value = "'0.003465"
$ws.Range("E48").Value = "'-1.18%"
$ws.Range("D49").Value = "'0.003530"
$ws.Range("D50").Value = "'0.00002100"
$ws.Range("E50").Value = "'-0.16%"
$ws.Range("D51").Value = "'0.0002000"
$ws.Range("E51").Value = "'-0.16%"
